# Weekly update: insert 3 new price rows (one per Plátano quality/variety
# combo) for the latest reporting date, pushing the existing history down.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows right before the existing data that needs to shift
# down (old row 327 -> new row 330, ... old row 415 -> new row 418).
$ws.Rows("327:329").Insert()

# Columns that are constant for every record in this sheet (single
# market / product combination).
$mercadoId = 4
$mercado   = "Feria Lagunitas de Puerto Montt"
$region    = "Los Lagos"
$codreg    = 10
$tipo      = "Fruta"
$productoId = 100108
$producto   = "Tropicales y subtropicales"
$categoriaId = 100108006
$categoria   = "Plátano"
$unidad    = "$/caja 20 kilos"
$origen    = "Ecuador"
$kgUnidad  = 20

$fecha = Get-Date -Year 2022 -Month 3 -Day 22 -Hour 0 -Minute 0 -Second 0

function Set-Fila {
    param(
        [int]$Fila,
        [string]$Variedad,
        [string]$Calidad,
        [double]$Volumen,
        [double]$PrecioMin,
        [double]$PrecioMax,
        [double]$PrecioProm,
        [double]$PrecioKg
    )

    $ws.Cells.Item($Fila, 1).Value  = $mercadoId
    $ws.Cells.Item($Fila, 2).Value  = $mercado
    $ws.Cells.Item($Fila, 3).Value  = $region
    $ws.Cells.Item($Fila, 4).Value  = $fecha
    $ws.Cells.Item($Fila, 5).Value  = $codreg
    $ws.Cells.Item($Fila, 6).Value  = $tipo
    $ws.Cells.Item($Fila, 7).Value  = $productoId
    $ws.Cells.Item($Fila, 8).Value  = $producto
    $ws.Cells.Item($Fila, 9).Value  = $categoriaId
    $ws.Cells.Item($Fila, 10).Value = $categoria
    $ws.Cells.Item($Fila, 11).Value = $Variedad
    $ws.Cells.Item($Fila, 12).Value = $Calidad
    $ws.Cells.Item($Fila, 13).Value = $Volumen
    $ws.Cells.Item($Fila, 14).Value = $PrecioMin
    $ws.Cells.Item($Fila, 15).Value = $PrecioMax
    $ws.Cells.Item($Fila, 16).Value = $PrecioProm
    $ws.Cells.Item($Fila, 17).Value = $unidad
    $ws.Cells.Item($Fila, 18).Value = $origen
    $ws.Cells.Item($Fila, 19).Value = $PrecioKg
    $ws.Cells.Item($Fila, 20).Value = $kgUnidad
}

Set-Fila 327 "Barraganete"     "Primera"        300  25000 26000 25500 1275
Set-Fila 328 "Sin especificar" "Pintón"         400  19000 19000 19000 950
Set-Fila 329 "Sin especificar" "Primera Pintón" 1200 21000 22000 21500 1075
